# Apply "fixing scenario to test commissioning scenarios" edits
# to the "Coupling Parameters" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2040 -> 2045
$ws.Range("B4").Value = 2045

# Row 18: label text changes + flag flips FALSE -> TRUE
$ws.Range("A18").Value = "fix_demand_and_profiles_to_initial_year"
$ws.Range("B18").Value = $true

# Row 19: flag flips TRUE -> FALSE
$ws.Range("B19").Value = $false

# Row 10: start_tick_dismantling value 3 -> 50
$ws.Range("B10").Value = 50

# Restore default view (top-left cell back to A1) and select B4
$ws.Activate()
$ws.Range("B4").Select()
